# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with newly scraped quotes. Price cells are set via NumberFormat "@" (text)
# then reset to the "Normal" style so numeric-looking strings (e.g. "1.010",
# "0.07129") survive as text instead of being coerced into floats/doubles,
# matching the original inline-string cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.085.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.769.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3759"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3393"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07344"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.349"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.776.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.960"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001077"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.503"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.097.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.475"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.482"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.977.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.063"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.887"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.654"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.335"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6723"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06272"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2158"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.642"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.11%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6256"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.821"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.092"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.87%  "
